# Add a "hashcode" column (I) with a per-speaker hashcode/slug for each
# row of the Public Works speaker table.
#
# Column I, row 1 is the new header ("hashcode"); rows 2-27 hold one
# hashcode per speaker-row, matching the existing 26 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# header + per-row hashcode values, in row order (row 1 = header)
$hashcodes = @(
    "hashcode",
    "sethkoproski",
    "zachulibarri",
    "evesnyder",
    "jamesnagy",
    "benfried",
    "ligiacoelho",
    "meganbarrington",
    "sethstrickland",
    "danielasamur",
    "lukekeller",
    "hunteradams",
    "michaelcaporizzo",
    "mollyryan",
    "karlsmolenski",
    "daisyrosas",
    "alisonritterhaus",
    "zachulibarri2",
    "adamhawkins",
    "jayleeming",
    "danielsprocket",
    "heatherhuson",
    "alitahoward",
    "jimmyjordan",
    "adamhowell",
    "madisonfitzpatrick",
    "andresmontealegre"
)

for ($i = 0; $i -lt $hashcodes.Length; $i++) {
    $row = $i + 1
    $cell = $ws.Cells.Item($row, 9)
    $cell.Value = $hashcodes[$i]

    # Rows 1-23 carry a distinct (but format-equivalent) cell style versus
    # the plain default used by the rest of the sheet - the last 4 rows
    # (24-27) were added without it, matching the source workbook.
    if ($row -le 23) {
        $cell.ShrinkToFit = $false
    }
}

# Column D ("speaker") is a bit wider in the edited workbook.
$ws.Columns.Item(4).ColumnWidth = 24.85546875
